$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 611, shifting existing rows 611..638 down to 612..639.
$ws.Rows.Item(611).Insert()

# Populate the newly inserted row 611 with the new weekly price record.
$ws.Cells.Item(611, 1).Value = 4
$ws.Cells.Item(611, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(611, 3).Value = "Los Lagos"
$ws.Cells.Item(611, 4).Value = 45267
$ws.Cells.Item(611, 5).Value = 10
$ws.Cells.Item(611, 6).Value = 100112023
$ws.Cells.Item(611, 7).Value = "Brócoli"
$ws.Cells.Item(611, 8).Value = "Sin especificar"
$ws.Cells.Item(611, 9).Value = "Primera"
$ws.Cells.Item(611, 10).Value = 500
$ws.Cells.Item(611, 11).Value = 1600
$ws.Cells.Item(611, 12).Value = 1600
$ws.Cells.Item(611, 13).Value = 1600
$ws.Cells.Item(611, 14).Value = "`$/unidad"
$ws.Cells.Item(611, 15).Value = "Región Metropolitana"
$ws.Cells.Item(611, 16).Value = 1600
$ws.Cells.Item(611, 17).Value = 1
$ws.Cells.Item(611, 18).Value = "Hortaliza"
